$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 10002663  # was 11113821
$ws.Range("I86").Value = 14288312  # was 16669248
$ws.Range("J86").Value = 2816.3333  # was 2966.3333
$ws.Range("K86").Value = 14288312  # was 16669248
$ws.Range("L86").Value = 2816.3333  # was 2966.3333
$ws.Range("M86").Value = -14287189  # was -16668125
$ws.Range("N86").Value = -5062.3333  # was -5212.3333
$ws.Range("H89").Value = 10002663  # was 11113821
$ws.Range("I89").Value = 14288312  # was 16669248
$ws.Range("J89").Value = 2816.3333  # was 2966.3333
$ws.Range("K89").Value = 71441560  # was 83346240
$ws.Range("L89").Value = 14081.6665  # was 14831.6665
$ws.Range("M89").Value = -71435944  # was -83340624
$ws.Range("N89").Value = -25313.6665  # was -26063.6665
$ws.Range("H96").Value = 334746.84  # was 500863
$ws.Range("I96").Value = 500805.75  # was 667474.3
$ws.Range("J96").Value = 2629  # was 1029
$ws.Range("K96").Value = 1502417.25  # was 2002422.9
$ws.Range("L96").Value = 7887  # was 3087
$ws.Range("M96").Value = -1501044.25  # was -2001049.9
$ws.Range("N96").Value = -10633  # was -5833
$ws.Range("H106").Value = 27996  # was 100000
$ws.Range("I106").Value = 2994.5  # was 0
$ws.Range("J106").Value = 52997.5  # was 100000
$ws.Range("K106").Value = 2994.5  # was 0
$ws.Range("L106").Value = 52997.5  # was 100000
$ws.Range("M106").Value = -2363.5  # new cell
$ws.Range("N106").Value = -54259.5  # was -101262
$ws.Range("H113").Value = 4636.3687  # was 4480.524
$ws.Range("I113").Value = 2200  # was 2733.3333
$ws.Range("K113").Value = 2200  # was 2733.3333
$ws.Range("M113").Value = 1054  # was 520.6667000000002
$ws.Range("H131").Value = 386.55554  # was 334.125
$ws.Range("I131").Value = 422.5  # was 340
$ws.Range("J131").Value = 99  # was 293
$ws.Range("K131").Value = 1267.5  # was 1020
$ws.Range("L131").Value = 297  # was 879
$ws.Range("M131").Value = 3772.5  # was 4020
$ws.Range("N131").Value = -10377  # was -10959
$ws.Range("H132").Value = 2913.4146  # was 2989.725
$ws.Range("I132").Value = 1054.0605  # was 1091.3438
$ws.Range("K132").Value = 3162.1815  # was 3274.0314
$ws.Range("M132").Value = -632.1815000000001  # was -744.0314000000003
$ws.Range("H138").Value = 2697.2144  # was 2739.6155
$ws.Range("J138").Value = 3210  # was 3328.2222
$ws.Range("L138").Value = 9630  # was 9984.6666
$ws.Range("N138").Value = -19910  # was -20264.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4127.75  # was 4366.3335
$ws.Range("I61").Value = 4127.75  # was 4366.3335
$ws.Range("K61").Value = 4127.75  # was 4366.3335
$ws.Range("M61").Value = -3915.75  # was -4154.3335
$ws.Range("H102").Value = 3749.4  # was 4797
$ws.Range("I102").Value = 3888.2222  # was 4797
$ws.Range("J102").Value = 2500  # was 0
$ws.Range("K102").Value = 3888.2222  # was 4797
$ws.Range("L102").Value = 2500  # was 0
$ws.Range("M102").Value = -2266.2222  # was -3175
$ws.Range("N102").Value = -5744  # new cell
$ws.Range("H132").Value = 3970.5  # was 4335.6665
$ws.Range("I132").Value = 2875  # was 0
$ws.Range("K132").Value = 8625  # was 0
$ws.Range("M132").Value = -6095  # new cell
$ws.Range("H136").Value = 4127.75  # was 4366.3335
$ws.Range("I136").Value = 4127.75  # was 4366.3335
$ws.Range("K136").Value = 12383.25  # was 13099.0005
$ws.Range("M136").Value = -9833.25  # was -10549.0005
$ws.Range("H139").Value = 99997  # was 99998
$ws.Range("J139").Value = 99997  # was 99998
$ws.Range("L139").Value = 99997  # was 99998
$ws.Range("N139").Value = -110277  # was -110278

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1711.5  # was 1655.875
$ws.Range("I86").Value = 1434.8182  # was 1435.5454
$ws.Range("J86").Value = 2320.2  # was 2140.6
$ws.Range("K86").Value = 1434.8182  # was 1435.5454
$ws.Range("L86").Value = 2320.2  # was 2140.6
$ws.Range("M86").Value = -311.8181999999999  # was -312.5454
$ws.Range("N86").Value = -4566.2  # was -4386.6
$ws.Range("H89").Value = 1711.5  # was 1655.875
$ws.Range("I89").Value = 1434.8182  # was 1435.5454
$ws.Range("J89").Value = 2320.2  # was 2140.6
$ws.Range("K89").Value = 7174.090999999999  # was 7177.727
$ws.Range("L89").Value = 11601  # was 10703
$ws.Range("M89").Value = -1558.090999999999  # was -1561.727
$ws.Range("N89").Value = -22833  # was -21935
$ws.Range("H99").Value = 3072.3157  # was 3430.5
$ws.Range("I99").Value = 3237.4666  # was 3430.5
$ws.Range("J99").Value = 2453  # was 0
$ws.Range("K99").Value = 3237.4666  # was 3430.5
$ws.Range("L99").Value = 2453  # was 0
$ws.Range("M99").Value = -1739.4666  # was -1932.5
$ws.Range("N99").Value = -5449  # new cell
$ws.Range("H116").Value = 0  # was 80000
$ws.Range("J116").Value = 0  # was 80000
$ws.Range("L116").Value = 0  # was 80000
$ws.Range("N116").ClearContents()  # was -89178
$ws.Range("H134").Value = 2682.9333  # was 2929.1667
$ws.Range("I134").Value = 2403.3845  # was 2615
$ws.Range("K134").Value = 7210.1535  # was 7845
$ws.Range("M134").Value = -4675.1535  # was -5310

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6968.8  # was 8837
$ws.Range("I16").Value = 6610.6665  # was 8249.5
$ws.Range("J16").Value = 7506  # was 10012
$ws.Range("K16").Value = 6610.6665  # was 8249.5
$ws.Range("L16").Value = 7506  # was 10012
$ws.Range("M16").Value = -6323.6665  # was -7962.5
$ws.Range("N16").Value = -8080  # was -10586
$ws.Range("H26").Value = 10000  # was 8500
$ws.Range("J26").Value = 10000  # was 8500
$ws.Range("L26").Value = 10000  # was 8500
$ws.Range("N26").Value = -10574  # was -9074
$ws.Range("H58").Value = 5844.909  # was 6094.476
$ws.Range("I58").Value = 6101.9  # was 6391.263
$ws.Range("K58").Value = 6101.9  # was 6391.263
$ws.Range("M58").Value = -5898.9  # was -6188.263
$ws.Range("H59").Value = 17594  # was 16441.111
$ws.Range("J59").Value = 18742.5  # was 16871.25
$ws.Range("L59").Value = 18742.5  # was 16871.25
$ws.Range("N59").Value = -21032.5  # was -19161.25
$ws.Range("H86").Value = 8068.2856  # was 8158.154
$ws.Range("I86").Value = 7626.8  # was 7707.5557
$ws.Range("K86").Value = 7626.8  # was 7707.5557
$ws.Range("M86").Value = -6503.8  # was -6584.5557
$ws.Range("H89").Value = 8068.2856  # was 8158.154
$ws.Range("I89").Value = 7626.8  # was 7707.5557
$ws.Range("K89").Value = 38134  # was 38537.7785
$ws.Range("M89").Value = -32518  # was -32921.7785
$ws.Range("H105").Value = 3502.611  # was 3579.8572
$ws.Range("I105").Value = 3866.4211  # was 4036.8333
$ws.Range("K105").Value = 3866.4211  # was 4036.8333
$ws.Range("M105").Value = -2119.4211  # was -2289.8333
$ws.Range("H107").Value = 742.75  # was 770.7895
$ws.Range("I107").Value = 446.26666  # was 463.14285
$ws.Range("K107").Value = 446.26666  # was 463.14285
$ws.Range("M107").Value = 1473.73334  # was 1456.85715
$ws.Range("H113").Value = 6968.8  # was 8837
$ws.Range("I113").Value = 6610.6665  # was 8249.5
$ws.Range("J113").Value = 7506  # was 10012
$ws.Range("K113").Value = 6610.6665  # was 8249.5
$ws.Range("L113").Value = 7506  # was 10012
$ws.Range("M113").Value = -4440.6665  # was -6079.5
$ws.Range("N113").Value = -11846  # was -14352
$ws.Range("H136").Value = 5844.909  # was 6094.476
$ws.Range("I136").Value = 6101.9  # was 6391.263
$ws.Range("K136").Value = 18305.7  # was 19173.789
$ws.Range("M136").Value = -15755.7  # was -16623.789

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1030.75  # was 951.7692
$ws.Range("J107").Value = 936.9  # was 852.0909
$ws.Range("L107").Value = 2810.7  # was 2556.2727
$ws.Range("N107").Value = -6650.7  # was -6396.2727
$ws.Range("H114").Value = 9091870  # was 10001009
$ws.Range("I114").Value = 22223076  # was 25000900
$ws.Range("J114").Value = 1035.1538  # was 1081.8334
$ws.Range("K114").Value = 66669228  # was 75002700
$ws.Range("L114").Value = 3105.4614  # was 3245.5002
$ws.Range("M114").Value = -66665974  # was -74999446
$ws.Range("N114").Value = -9613.4614  # was -9753.5002
$ws.Range("H121").Value = 20835344  # was 20835388
$ws.Range("J121").Value = 2785.6365  # was 2849.0908
$ws.Range("L121").Value = 8356.9095  # was 8547.2724
$ws.Range("N121").Value = -10976.9095  # was -11167.2724

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 8751.75  # was 14877.25
$ws.Range("I55").Value = 10002.5  # was 10169.833
$ws.Range("J55").Value = 4999.5  # was 28999.5
$ws.Range("K55").Value = 10002.5  # was 10169.833
$ws.Range("L55").Value = 4999.5  # was 28999.5
$ws.Range("M55").Value = -9675.5  # was -9842.833000000001
$ws.Range("N55").Value = -5653.5  # was -29653.5
$ws.Range("H134").Value = 90000  # was 0
$ws.Range("J134").Value = 90000  # was 0
$ws.Range("L134").Value = 270000  # was 0
$ws.Range("N134").Value = -275070  # new cell

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1566.7858  # was 1522.2667
$ws.Range("I22").Value = 2071.1667  # was 1903.7142
$ws.Range("K22").Value = 2071.1667  # was 1903.7142
$ws.Range("M22").Value = -1776.1667  # was -1608.7142
$ws.Range("H27").Value = 1566.7858  # was 1522.2667
$ws.Range("I27").Value = 2071.1667  # was 1903.7142
$ws.Range("K27").Value = 2071.1667  # was 1903.7142
$ws.Range("M27").Value = -1964.1667  # was -1796.7142
$ws.Range("H50").Value = 422495  # was 42495
$ws.Range("J50").Value = 422495  # was 42495
$ws.Range("L50").Value = 422495  # was 42495
$ws.Range("N50").Value = -423769  # was -43769
$ws.Range("H132").Value = 4742.7104  # was 4846.7837
$ws.Range("J132").Value = 3731.1667  # was 4299
$ws.Range("L132").Value = 11193.5001  # was 12897
$ws.Range("N132").Value = -16253.5001  # was -17957

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3732.5715  # was 3995
$ws.Range("I81").Value = 3725.8  # was 3995
$ws.Range("J81").Value = 3749.5  # was 3995
$ws.Range("K81").Value = 7451.6  # was 7990
$ws.Range("L81").Value = 7499  # was 7990
$ws.Range("M81").Value = -6390.6  # was -6929
$ws.Range("N81").Value = -9621  # was -10112
$ws.Range("H84").Value = 3732.5715  # was 3995
$ws.Range("I84").Value = 3725.8  # was 3995
$ws.Range("J84").Value = 3749.5  # was 3995
$ws.Range("K84").Value = 37258  # was 39950
$ws.Range("L84").Value = 37495  # was 39950
$ws.Range("M84").Value = -31954  # was -34646
$ws.Range("N84").Value = -48103  # was -50558
